$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells whose new price values would otherwise be auto-parsed as
# numbers by Excel stay formatted/stored as plain text (matching the
# original inline-string cells).
$textCells = @("D5", "D6", "D8", "D14", "D20", "D21", "D25", "D26", "D27", "D28", "D31", "D35", "D38", "D39", "D45", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Price (column D) updates
$ws.Range("D2").Value = "63.020.69"
$ws.Range("D3").Value = "2.459.69"
$ws.Range("D5").Value = "576.35"
$ws.Range("D6").Value = "146.58"
$ws.Range("D8").Value = "0.541"
$ws.Range("D9").Value = "2.459.15"
$ws.Range("D14").Value = "28.99"
$ws.Range("D16").Value = "2.906.52"
$ws.Range("D17").Value = "62.922.90"
$ws.Range("D18").Value = "2.459.23"
$ws.Range("D20").Value = "11.10"
$ws.Range("D21").Value = "330.29"
$ws.Range("D25").Value = "66.46"
$ws.Range("D26").Value = "665.58"
$ws.Range("D27").Value = "1.17"
$ws.Range("D28").Value = "8.96"
$ws.Range("D29").Value = "0.0₃0999"
$ws.Range("D30").Value = "2.580.40"
$ws.Range("D31").Value = "8.18"
$ws.Range("D35").Value = "1.54"
$ws.Range("D38").Value = "5.51"
$ws.Range("D39").Value = "153.38"
$ws.Range("D42").Value = "0.0₆0349"
$ws.Range("D45").Value = "42.29"
$ws.Range("D48").Value = "146.60"
$ws.Range("D50").Value = "20.72"

# Volume/1h change (column E) updates
$ws.Range("E2").Value = "  +2.48%  "
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("E6").Value = "  +2.33%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").Value = "  +1.57%  "
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  +1.17%  "
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("E14").Value = "  +9.63%  "
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("E16").Value = "  +2.23%  "
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E26").Value = "  +8.38%  "
$ws.Range("E27").Value = "  +17.44%  "
$ws.Range("E28").Value = "  +8.16%  "
$ws.Range("E29").Value = "  +4.79%  "
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("E31").Value = "  +1.91%  "
$ws.Range("E32").Value = "  +3.92%  "
$ws.Range("E33").Value = "  +5.10%  "
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("E35").Value = "  +4.60%  "
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("E38").Value = "  +2.91%  "
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("E42").Value = "  +23.40%  "
$ws.Range("E43").Value = "  +6.80%  "
$ws.Range("E44").Value = "  +4.14%  "
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +27.54%  "
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  +4.08%  "
$ws.Range("E51").Value = "  +2.08%  "

